$d = $word.ActiveDocument

# The document's final paragraph (right before the section break) is an
# empty paragraph whose mark already carries sz=24/szCs=24/lang=en-US
# formatting. Add a new run there with the "www.getbootstrap.com" source
# link, carrying that same run formatting.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Collapse to a zero-length insertion point right before the paragraph
# mark, so the new run lands inside this paragraph (not a new one).
$insertPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

# Use InsertXML so the run gets exact sz/szCs/lang run-properties (plain
# Range.Text + Font.* setters can't express szCs/lang).
$runXml = '<?xml version="1.0" standalone="yes"?>' +
          '<?mso-application progid="Word.Document"?>' +
          '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body><w:p><w:r>' +
          '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>' +
          '<w:t>www.getbootstrap.com</w:t>' +
          '</w:r></w:p></w:body></w:document>' +
          '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($runXml)
